$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Carlos"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "02/05/2005"

$ws.Range("D2").Value = "00:00"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 40
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "Pendente"
